$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 425
$ws.Range("I45").Value = 400
$ws.Range("J45").Value = 500
$ws.Range("K45").Value = 1200
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -1008
$ws.Range("N45").Value = -1884

$ws.Range("H80").Value = 3797.1765
$ws.Range("I80").Value = 2081.375
$ws.Range("J80").Value = 5322.3335
$ws.Range("K80").Value = 6244.125
$ws.Range("L80").Value = 15967.0005
$ws.Range("M80").Value = -5246.125
$ws.Range("N80").Value = -17963.0005

$ws.Range("H83").Value = 3797.1765
$ws.Range("I83").Value = 2081.375
$ws.Range("J83").Value = 5322.3335
$ws.Range("K83").Value = 18732.375
$ws.Range("L83").Value = 47901.0015
$ws.Range("M83").Value = -13740.375
$ws.Range("N83").Value = -57885.0015

$ws.Range("H88").Value = 1978.8
$ws.Range("I88").Value = 596
$ws.Range("J88").Value = 2324.5
$ws.Range("K88").Value = 596
$ws.Range("L88").Value = 2324.5
$ws.Range("M88").Value = -190
$ws.Range("N88").Value = -3136.5

$ws.Range("H91").Value = 1978.8
$ws.Range("I91").Value = 596
$ws.Range("J91").Value = 2324.5
$ws.Range("K91").Value = 596
$ws.Range("L91").Value = 2324.5
$ws.Range("M91").Value = 808
$ws.Range("N91").Value = -5132.5

$ws.Range("H100").Value = 3075
$ws.Range("I100").Value = 2650
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 2650
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -2109
$ws.Range("N100").Value = -4582

$ws.Range("H137").Value = 1867.6666
$ws.Range("I137").Value = 1814.75
$ws.Range("J137").Value = 1973.5
$ws.Range("K137").Value = 5444.25
$ws.Range("L137").Value = 5920.5
$ws.Range("M137").Value = -2894.25
$ws.Range("N137").Value = -11020.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5114.44
$ws.Range("I61").Value = 3868
$ws.Range("J61").Value = 9061.5
$ws.Range("K61").Value = 3868
$ws.Range("L61").Value = 9061.5
$ws.Range("M61").Value = -3656

$ws.Range("H110").Value = 12333970
$ws.Range("I110").Value = 12333970
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 12333970
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -12331925

$ws.Range("H122").Value = 2493.3333
$ws.Range("I122").Value = 980
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 2940
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -490
$ws.Range("N122").Value = -14650

$ws.Range("H136").Value = 5114.44
$ws.Range("I136").Value = 3868
$ws.Range("J136").Value = 9061.5
$ws.Range("K136").Value = 11604
$ws.Range("L136").Value = 27184.5
$ws.Range("M136").Value = -9054

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 319.8
$ws.Range("I80").Value = 285.55554
$ws.Range("J80").Value = 371.16666
$ws.Range("K80").Value = 285.55554
$ws.Range("L80").Value = 371.16666
$ws.Range("M80").Value = 712.4444599999999

$ws.Range("H83").Value = 319.8
$ws.Range("I83").Value = 285.55554
$ws.Range("J83").Value = 371.16666
$ws.Range("K83").Value = 1427.7777
$ws.Range("L83").Value = 1855.8333
$ws.Range("M83").Value = 3564.2223

$ws.Range("H94").Value = 1625.5714
$ws.Range("I94").Value = 1579.8334
$ws.Range("J94").Value = 1900
$ws.Range("K94").Value = 1579.8334
$ws.Range("L94").Value = 1900
$ws.Range("M94").Value = -1128.8334
$ws.Range("N94").Value = -2802

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2039.6316
$ws.Range("I31").Value = 1813.2307
$ws.Range("J31").Value = 2530.1667
$ws.Range("K31").Value = 1813.2307
$ws.Range("L31").Value = 2530.1667
$ws.Range("M31").Value = -1518.2307
$ws.Range("N31").Value = -3120.1667

$ws.Range("H34").Value = 2039.6316
$ws.Range("I34").Value = 1813.2307
$ws.Range("J34").Value = 2530.1667
$ws.Range("K34").Value = 1813.2307
$ws.Range("L34").Value = 2530.1667
$ws.Range("M34").Value = -1611.2307
$ws.Range("N34").Value = -2934.1667

$ws.Range("H105").Value = 2830.2693
$ws.Range("I105").Value = 2255.8572
$ws.Range("J105").Value = 3500.4167
$ws.Range("K105").Value = 2255.8572
$ws.Range("L105").Value = 3500.4167
$ws.Range("M105").Value = -508.8571999999999

$ws.Range("H122").Value = 2511.6667
$ws.Range("I122").Value = 2511.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7535.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5085.000100000001

$ws.Range("H125").Value = 140326
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 140326
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 140326
$ws.Range("N125").Value = -145246

$ws.Range("H132").Value = 1380.5151
$ws.Range("I132").Value = 1314.2812
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 3942.8436
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -1412.8436

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 792.5
$ws.Range("I14").Value = 792.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2377.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2204.5

$ws.Range("H23").Value = 355.14285
$ws.Range("I23").Value = 77.2
$ws.Range("J23").Value = 1050
$ws.Range("K23").Value = 231.6
$ws.Range("L23").Value = 3150
$ws.Range("M23").Value = 3.399999999999977
$ws.Range("N23").Value = -3620

$ws.Range("H37").Value = 59952
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 59952
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 179856
$ws.Range("N37").Value = -180080

$ws.Range("H132").Value = 4439.2
$ws.Range("I132").Value = 3398.8333
$ws.Range("J132").Value = 5999.75
$ws.Range("K132").Value = 30589.4997
$ws.Range("L132").Value = 53997.75
$ws.Range("M132").Value = -28059.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

$ws.Range("H68").Value = 112634
$ws.Range("I68").Value = 100268
$ws.Range("J68").Value = 125000
$ws.Range("K68").Value = 100268
$ws.Range("L68").Value = 125000
$ws.Range("M68").Value = -99457
$ws.Range("N68").Value = -126622

$ws.Range("H70").Value = 7320
$ws.Range("I70").Value = 6729
$ws.Range("J70").Value = 7714
$ws.Range("K70").Value = 6729
$ws.Range("L70").Value = 7714
$ws.Range("M70").Value = -6459
$ws.Range("N70").Value = -8254

$ws.Range("H71").Value = 112634
$ws.Range("I71").Value = 100268
$ws.Range("J71").Value = 125000
$ws.Range("K71").Value = 300804
$ws.Range("L71").Value = 375000
$ws.Range("M71").Value = -296748
$ws.Range("N71").Value = -383112

$ws.Range("H73").Value = 7320
$ws.Range("I73").Value = 6729
$ws.Range("J73").Value = 7714
$ws.Range("K73").Value = 6729
$ws.Range("L73").Value = 7714
$ws.Range("M73").Value = -5793
$ws.Range("N73").Value = -9586

$ws.Range("H92").Value = 23083.666
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 23083.666
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 23083.666
$ws.Range("N92").Value = -26827.666

$ws.Range("H113").Value = 300
$ws.Range("I113").Value = 300
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 300
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1870

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1513.3846
$ws.Range("I16").Value = 1472.8334
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1472.8334
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1302.8334

$ws.Range("H55").Value = 213.5
$ws.Range("I55").Value = 162.4
$ws.Range("J55").Value = 298.66666
$ws.Range("K55").Value = 162.4
$ws.Range("L55").Value = 298.66666
$ws.Range("M55").Value = 10.59999999999999
$ws.Range("N55").Value = -644.66666

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H82").Value = 1586.8235
$ws.Range("I82").Value = 1591.2858
$ws.Range("J82").Value = 1566
$ws.Range("K82").Value = 1591.2858
$ws.Range("L82").Value = 1566
$ws.Range("M82").Value = -1230.2858

$ws.Range("H85").Value = 1586.8235
$ws.Range("I85").Value = 1591.2858
$ws.Range("J85").Value = 1566
$ws.Range("K85").Value = 1591.2858
$ws.Range("L85").Value = 1566
$ws.Range("M85").Value = -343.2858000000001

$ws.Range("H122").Value = 9100.691999999999
$ws.Range("I122").Value = 12136.333
$ws.Range("J122").Value = 6498.7144
$ws.Range("K122").Value = 36408.999
$ws.Range("L122").Value = 19496.1432
$ws.Range("M122").Value = -33958.999
$ws.Range("N122").Value = -24396.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 11173
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 11173
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 11173
$ws.Range("N76").Value = -11803

$ws.Range("H79").Value = 11173
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 11173
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 11173
$ws.Range("N79").Value = -13357
